$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B1 currently holds the number 15; it should become the text "15" instead,
# while keeping its existing cell style untouched. A direct
# "$ws.Range('B1').Value = '15'" would get auto-coerced back to a number by
# Excel's smart typing, and forcing text via NumberFormat/quote-prefix would
# stamp a brand-new style on the cell. Copying the existing text cell C8
# (which already stores the text "15") and pasting values-only reproduces
# the same text value/type on B1 without touching its style.
$ws.Range("C8").Copy()
$ws.Range("B1").PasteSpecial(-4163)

# Add the new product row 9: "Blem", 10, 50 (quantity/price columns stored
# as text, matching the existing data rows 4-8). Re-use B4/C4, which already
# hold the text values "10" and "50", as paste sources so the new cells end
# up with the same (un-styled) text type as the rest of the table.
$ws.Range("A9").Value = "Blem"

$ws.Range("B4").Copy()
$ws.Range("B9").PasteSpecial(-4163)

$ws.Range("C4").Copy()
$ws.Range("C9").PasteSpecial(-4163)

$excel.CutCopyMode = $false
